$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Summer of Magic (PSUM)"
$ws.Range("A2").Value = "Demonic Tutor"
$ws.Range("A3").Value = "Goblin Piledriver"
$ws.Range("A4").Value = "Living Wish"
$ws.Range("A5").Value = "Mind's Desire"
$ws.Range("A6").Value = "Orim's Chant"
